# Apple DCF Tiny.xlsx - "Add files via upload" edit
#
# The underlying change is: the two "downside scenario" checkboxes on the
# Assumptions sheet (linked to cells $G$24 and $G$25) get checked/turned on.
# That ripples, through the existing formulas, into Revenue Model,
# Financials, FCF, WACC, "FV Calculation", and finally the Dashboard summary
# and its chart. We also re-select the Dashboard tab as the active sheet
# (it was the Assumptions tab that was active/selected before).

$wb = $excel.ActiveWorkbook

$wsAssumptions = $wb.Worksheets.Item("Assumptions")

# Check the two form-control checkboxes (Assumptions!G24 / Assumptions!G25)
# by setting their linked cells to TRUE -- this is what Excel does when a
# linked CheckBox is clicked/checked.
$wsAssumptions.Range("G24").Value = $true
$wsAssumptions.Range("G25").Value = $true

# Also flip the checkbox form controls themselves so their visual state
# agrees with the linked cells.
foreach ($shp in $wsAssumptions.Shapes) {
    if ($shp.Name -eq "Check Box 12" -or $shp.Name -eq "Check Box 13") {
        $shp.ControlFormat.Value = 1
    }
}

# Recalculate everything so dependent sheets/chart pick up the new values.
$excel.CalculateFull()

# Refresh the "FCF Fair Values" chart on the FV Calculation sheet so its
# cached series data matches the recalculated cells.
$wsFV = $wb.Worksheets.Item("FV Calculation")
foreach ($co in $wsFV.ChartObjects()) {
    $co.Chart.Refresh()
}

# The active/selected sheet moves from Assumptions back to Dashboard.
$wsDashboard = $wb.Worksheets.Item("Dashboard")
$wsDashboard.Activate()
